$d = $word.ActiveDocument

# 1. "Oats 3Kg" -> "Quaker Oats 3Kg"
$d.Paragraphs.Item(4).Range.Find.Execute("Oats 3Kg", $true, $false, $false, $false, $false, $true, 1, $false, "Quaker Oats 3Kg", 2)

# 2. "Approximately $20 per 1kg" -> "$26"
$d.Paragraphs.Item(6).Range.Find.Execute("20", $true, $false, $false, $false, $false, $true, 1, $false, "26", 2)

# 3. "Price: $60" (Oats) -> "$78"
$d.Paragraphs.Item(7).Range.Find.Execute("60", $true, $false, $false, $false, $false, $true, 1, $false, "78", 2)

# 4. "Approximately ~ $10 per Piece" -> "$12"
$d.Paragraphs.Item(12).Range.Find.Execute("10 ", $true, $false, $false, $false, $false, $true, 1, $false, "12 ", 2)

# 5. "Price: $60" (Chicken) -> "$72"
$d.Paragraphs.Item(13).Range.Find.Execute("60", $true, $false, $false, $false, $false, $true, 1, $false, "72", 2)

# 6. "Approximate $11 per can" -> "$9 per can"
$d.Paragraphs.Item(22).Range.Find.Execute("Approximate $11 per can", $true, $false, $false, $false, $false, $true, 1, $false, "Approximate $9 per can", 2)

# 7. "Price: $33" -> "$27"
$d.Paragraphs.Item(23).Range.Find.Execute("Price: $33", $true, $false, $false, $false, $false, $true, 1, $false, "Price: $27", 2)

# 8. "Approximate $13 per Can" -> "$14 per Can"
$d.Paragraphs.Item(30).Range.Find.Execute("Approximate $13 per Can", $true, $false, $false, $false, $false, $true, 1, $false, "Approximate $14 per Can", 2)

# 9. "Price $39" -> "Price $42"
$d.Paragraphs.Item(31).Range.Find.Execute("Price $39", $true, $false, $false, $false, $false, $true, 1, $false, "Price $42", 2)

# 10. Total price "209.9" -> "236.9"
$d.Paragraphs.Item(42).Range.Find.Execute("209.9", $true, $false, $false, $false, $false, $true, 1, $false, "236.9", 2)

# 11. "Note: Some of the phot|os are Snap from the Huawei Mobile Phone" (two runs, bookmark
#     in between) -> single merged run, bookmark removed from here.
$p11 = $d.Paragraphs.Item(45)
$p11Start = $p11.Range.Start
$r11Start = $p11Start + 6
$r11End = $p11.Range.End - 1
$r11 = $d.Range($r11Start, $r11End)
$r11.Text = "TEMP_PLACEHOLDER_TEXT"
$p11b = $d.Paragraphs.Item(45)
$p11bStart = $p11b.Range.Start
$r11b = $d.Range($p11bStart + 6, $p11bStart + 6 + 21)
$r11b.Text = "Some of the photos are Snap from the Huawei Mobile Phone"

# 12. "Updated:" -> split into "Upd" + bookmark(_GoBack) + "ated:"
$p12 = $d.Paragraphs.Item(43)
$p12Start = $p12.Range.Start
$bmPos = $p12Start + 3
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
